$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.084.19"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "1.556.71"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "'292.13"
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("D7").Value = "'0.3980"
$ws.Range("E7").Value = "  +5.31%  "
$ws.Range("D8").Value = "'0.3229"
$ws.Range("E8").Value = "  -1.98%  "
$ws.Range("D9").Value = "'43.91"
$ws.Range("E9").Value = "  -2.12%  "
$ws.Range("D10").Value = "'0.07313"
$ws.Range("E10").Value = "  -1.23%  "
$ws.Range("D11").Value = "'1.083"
$ws.Range("E11").Value = "  -5.55%  "
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Value = "'18.91"
$ws.Range("E13").Value = "  -7.42%  "
$ws.Range("D14").Value = "'5.682"
$ws.Range("E14").Value = "  -3.51%  "
$ws.Range("D15").Value = "'0.00001137"
$ws.Range("E15").Value = "  +5.24%  "
$ws.Range("D16").Value = "'6.647"
$ws.Range("E16").Value = "  -2.03%  "
$ws.Range("D17").Value = "1.552.72"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("D19").Value = "'83.78"
$ws.Range("E19").Value = "  -3.15%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "'6.304"
$ws.Range("E21").Value = "  -2.04%  "
$ws.Range("D22").Value = "'15.72"
$ws.Range("E22").Value = "  -3.51%  "
$ws.Range("E23").Value = "  -3.94%  "
$ws.Range("D24").Value = "22.095.86"
$ws.Range("D25").Value = "'2.364"
$ws.Range("E25").Value = "  +2.81%  "
$ws.Range("E26").Value = "  -6.03%  "
$ws.Range("D27").Value = "'148.52"
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("D28").Value = "'18.63"
$ws.Range("E28").Value = "  -3.79%  "
$ws.Range("D29").Value = "'4.897"
$ws.Range("E29").Value = "  -0.72%  "
$ws.Range("D30").Value = "1.731.44"
$ws.Range("E30").Value = "  +1.60%  "
$ws.Range("D31").Value = "'119.13"
$ws.Range("E31").Value = "  -3.35%  "
$ws.Range("D32").Value = "'1.013"
$ws.Range("E32").Value = "  -6.63%  "
$ws.Range("D33").Value = "'5.794"
$ws.Range("E33").Value = "  -2.85%  "
$ws.Range("D34").Value = "'0.08351"
$ws.Range("E34").Value = "  +1.48%  "
$ws.Range("D35").Value = "'1.624"
$ws.Range("E35").Value = "  -15.89%  "
$ws.Range("D36").Value = "'9.099"
$ws.Range("E36").Value = "  -3.90%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.06130"
$ws.Range("E37").Value = "  -3.61%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02271"
$ws.Range("E38").Value = "  -4.32%  "
$ws.Range("D39").Value = "'5.142"
$ws.Range("E39").Value = "  -4.27%  "
$ws.Range("D40").Value = "'1.218"
$ws.Range("E40").Value = "  -2.69%  "
$ws.Range("D41").Value = "'0.2058"
$ws.Range("E41").Value = "  -4.97%  "
$ws.Range("D42").Value = "'1.000"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").Value = "'10.76"
$ws.Range("E43").Value = "  -2.94%  "
$ws.Range("D44").Value = "'0.5851"
$ws.Range("E44").Value = "  -4.27%  "
$ws.Range("D45").Value = "'3.760"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "'13.05"
$ws.Range("E46").Value = "  -6.20%  "
$ws.Range("D47").Value = "'0.5600"
$ws.Range("E47").Value = "  -5.80%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.917"
$ws.Range("E48").Value = "  -4.18%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'118.99"
$ws.Range("E49").Value = "  -3.40%  "
$ws.Range("E50").Value = "  -3.55%  "
$ws.Range("D51").Value = "'0.06848"
$ws.Range("E51").Value = "  -3.74%  "
